# Translate the Swahili facilitator-guide table labels (and a couple of
# body strings) into English, per the supplied diff.
#
# Find.Execute signature used below:
#   Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,
#           MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)
# Wrap = 1 (wdFindContinue), Replace = 2 (wdReplaceAll)

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "Kichwa cha Video" "Video Title"
Replace-Text "Mada" "Topic"
Replace-Text "Malengo" "Aim(s)"
Replace-Text "Urefu" "Length"
Replace-Text "Mahali pa Kambi" "Camp Location"
Replace-Text "Wawezeshaji" "Facilitators"
Replace-Text "N. ya wanafunzi" "N. of students"
Replace-Text "Tarehe" "Date"
Replace-Text "Rasilimali" "Resources"
Replace-Text "inahitajika" "needed"
Replace-Text "Maandalizi" "Preparations"
Replace-Text "Muda wa video" "Video time"
Replace-Text "Mwezeshaji anafanya nini" "What facilitator does"
Replace-Text "Wanachofanya wanafunzi" "What learners do"
Replace-Text "Utangulizi Mkuu wa Video ya VMC" "General VMC Video Introduction"
Replace-Text "Utangulizi wa Video" "Video Introduction"
Replace-Text "Kusaidia mchakato, kuchochea mawazo" "Assist the process, provoke thoughts"
Replace-Text "Suluhisho" "Solution"
